$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new "User" column between "CommitteeName" (A) and "Status" (B),
# shifting the old Status column to C.
$ws.Columns.Item(2).Insert()

$ws.Range("B1").Value = "User"
$ws.Range("B2").Value = "Harshita sharma"
$ws.Range("B3").Value = "ketan Sali"
$ws.Range("B4").Value = "Himanshu"
$ws.Range("B5").Value = "juku"
$ws.Range("B6").Value = "aaaaaa"
$ws.Range("B7").Value = "Akshay Bhagwat"
$ws.Range("B8").Value = "Manoj N"

$ws.Columns.Item(1).AutoFit()

$ws.Range("B9").Select()
